$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAD")

# Update "Used" hours for previous month (2022 June = col Q, 2022 July = col R)
$ws.Range("Q5").Value = 35
$ws.Range("R5").Value = 15

$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = 20

$ws.Range("Q7").Value = 135.234
$ws.Range("R7").Value = 30

$ws.Range("Q8").Value = 194.55
$ws.Range("R8").Value = 1.5

$ws.Range("Q9").Value = 150.5
$ws.Range("R9").Value = 36

# Update selection to match new active cell
[void]$ws.Range("Q6").Select()

# Update cell style name (localized "Normal" -> Czech "Normální")
$wb.Styles.Item("Normal").Name = "Normální"

